# Daily attendance processing - 2026-01-27 01:47:03
# Swap the order of "Recorded By" contributors from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" for every row in column G where that exact
# value appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
